$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the hourly crypto snapshot: Price (column D) and Volume(1h) % (column E).
# Column D holds text (e.g. "2.215.55" European thousands grouping), not real numbers,
# so each Price cell is forced to text format before the write, then the format is
# reset back to Normal so no stray number formatting lingers on the cell.

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "41.940.29"
$cell.Style = "Normal"
$ws.Range("E2").Value = "  -0.70%  "

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.214.56"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.42%  "

$ws.Range("E4").Value = "  +0.17%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "241.83"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -1.98%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "73.10"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "

$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  -2.05%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "42.31"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -0.65%  "

$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "0.0955"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  +0.52%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "7.05"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.71%  "

$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.103"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  +0.17%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.548.85"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "14.28"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.44%  "

$ws.Range("E16").Value = "  -1.93%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "2.213.17"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -2.75%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "41.855.21"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  -0.63%  "

$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.0000107"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +5.61%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "6.21"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.24%  "

$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "72.91"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +0.70%  "

$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "10.60"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +19.30%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "230.49"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "2.09"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -6.13%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "11.88"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.93%  "

$ws.Range("E26").Value = "  +0.07%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "3.68"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  +1.59%  "

$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "2.27"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("E29").Value = "  -2.99%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "168.20"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.78%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "20.49"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -1.89%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "5.63"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +7.02%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.0797"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -2.85%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "29.82"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -3.37%  "

$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("E36").Value = "  -9.83%  "

$ws.Range("E37").Value = "  -3.86%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.0302"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -4.25%  "

$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "13.75"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -0.61%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "65.94"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +5.08%  "

$ws.Range("E41").Value = "  -2.16%  "

$ws.Range("E42").Value = "  -2.25%  "

$ws.Range("E43").Value = "  -2.89%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "8.79"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +1.21%  "

$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "105.41"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -1.54%  "

$ws.Range("E46").Value = "  -2.11%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "2.44"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +5.50%  "

$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("E49").Value = "  -0.68%  "

$ws.Range("E50").Value = "  -0.12%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "2.423.25"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.44%  "
